# "Reiter und Text zugefügt"
# - Duplicate the existing sheet so the copy (placed before the original)
#   becomes "Buch" and keeps the original's underlying xl-part "flavor"
#   (this is what makes both resulting sheet parts come out with the same
#   shape/namespacing as the original single-sheet workbook had).
# - The original sheet (now pushed to 2nd position) becomes "Film".
# - Fill in the two sheets' cell content and fix up selection/active tab.
# - Resize the workbook window per the captured window geometry.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the current (only) sheet; places the copy right before $ws1.
$ws1.Copy($ws1) | Out-Null

# After the copy, position 1 is the new copy, position 2 is the original.
$buch = $wb.Worksheets.Item(1)
$film = $wb.Worksheets.Item(2)

$buch.Name = "Buch"
$film.Name = "Film"

# "Film" sheet just carries its own tab title as a label in A1.
$film.Range("A1").Value = "Film"

# "Buch" sheet: a tiny two-column listing.
$buch.Range("B2").Value = "Buch1"
$buch.Range("B3").Value = "Buch2"
$buch.Range("C1").Value = "Seite"
$buch.Range("D1").Value = "Datum"

# Make "Buch" the active/selected tab with C2 highlighted, matching the
# saved view state.
$buch.Select() | Out-Null
$buch.Range("C2").Select() | Out-Null

# Restore/resize the workbook window geometry.
$aw = $excel.ActiveWindow
$aw.Width = 28800
$aw.Height = 12300
